$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.497.20'
$ws.Range("E2").Value = '  +5.10%  '

# Row 3
$ws.Range("D3").Value = '1.724.40'
$ws.Range("E3").Value = '  +4.07%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5384'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.76%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2684'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.70%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06603'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.73'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07728'
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.641'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.32%  '

# Row 13
$ws.Range("D13").Value = '1.740.52'
$ws.Range("E13").Value = '  +4.90%  '

# Row 14
$ws.Range("D14").Value = '1.960.88'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5881'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.52%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8290'
$ws.Range("E16").Value = '  +1.03%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.73%  '

# Row 18
$ws.Range("D18").Value = '27.512.27'
$ws.Range("E18").Value = '  +5.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +15.00%  '

# Row 20
$ws.Range("E20").Value = '  +0.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.735'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.107'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.46%  '

# Row 24
$ws.Range("E24").Value = '  +0.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.65%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1233'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.67%  '

# Row 27
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.690'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.48%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.417'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.02%  '

# Row 29
$ws.Range("E29").Value = '  +4.35%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05555'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.12%  '

# Row 31
$ws.Range("E31").Value = '  +2.45%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.544'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.11%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.472'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.660'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.06%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9597'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.62%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.446'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.82%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.817'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.39%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5935'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.28%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01645'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.51%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.865'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8562'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.35%  '

# Row 42
$ws.Range("D42").Value = '1.055.99'
$ws.Range("E42").Value = '  +2.71%  '

# Row 43
$ws.Range("E43").Value = '  +0.02%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '

# Row 45
$ws.Range("D45").Value = '1.867.07'
$ws.Range("E45").Value = '  +3.98%  '

# Row 46
$ws.Range("E46").Value = '  +10.72%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '58.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.207'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.97%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4440'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.13%  '

# Row 50
$ws.Range("E50").Value = '  -0.14%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05272'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.42%  '
